$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the single ticket row (row 2) to the new "Sway" application ticket.
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = "Sway"
$ws.Range("E2").Value = "Microsoft Sway ended with an error is not able to open charts"
